# edit.ps1 - apply "small update for lab #2" changes
$p = $ppt.ActivePresentation

# =========================================================================
# Slide 1 ("Lab #2:" intro slide) -- shape 1 (TextBox 3)
# =========================================================================
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange

# --- merge the three leading runs of the "By the beginning..." paragraph
#     (paragraph 4) into a single run; keep the hyperlinked e-mail run
#     that follows untouched.
$para4 = $tr1.Paragraphs(4)
$leadRange = $para4.Characters(1, 65)
$leadRange.Text = "By the beginning of the next lab (Feb. 3), send what you have to "

# --- split the trailing ".." off of the "Send your code..." paragraph
#     (paragraph 6) into its own run.
$para6 = $tr1.Paragraphs(6)
$dotsRange = $tr1.Characters($para6.Start + 43, 2)
$dotsRange.Text = ".."

# --- append a new blank paragraph and a new paragraph reminding students
#     to put "Lab #2" in the e-mail subject line.
$para6 = $tr1.Paragraphs(6)
$makeSureText = "Make sure the text “Lab #2” is in the subject line… "
$para6.InsertAfter("`rZZZ_PLACEHOLDER_ZZZ`r$makeSureText") | Out-Null

$blankPara = $tr1.Paragraphs(7)
$blankPara.Text = ""

$makeSurePara = $tr1.Paragraphs(8)
$makeSureLeadRange = $tr1.Characters($makeSurePara.Start, 38)
$makeSureLeadRange.Text = "Make sure the text “Lab #2” is in the "

# =========================================================================
# Slide 3 (clinical-trial probability questions) -- shape 1 (TextBox 3)
# =========================================================================
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(1)
$tr3 = $sh3.TextFrame.TextRange

# Each of these paragraphs is currently split across 2-3 runs that carry
# identical formatting; collapse each back down to a single run. (First
# re-point at a placeholder string so the no-op/short-circuit on an
# unchanged Text value doesn't skip the underlying run merge.)
$tr3.Paragraphs(6).Text = "ZZZ_PLACEHOLDER_ZZZ"
$tr3.Paragraphs(6).Text = "`t(2A) Plot out the probability density function with the x-axis the number of"

$tr3.Paragraphs(7).Text = "ZZZ_PLACEHOLDER_ZZZ"
$tr3.Paragraphs(7).Text = "`tpatients that survive. "

$tr3.Paragraphs(9).Text = "ZZZ_PLACEHOLDER_ZZZ"
$tr3.Paragraphs(9).Text = "`t(2B) What is the p-value for a null hypothesis that the drug has no effect."

$tr3.Paragraphs(12).Text = "ZZZ_PLACEHOLDER_ZZZ"
$tr3.Paragraphs(12).Text = "`t(2C) What is the p-value for a null hypothesis that the drug does not "
